# Update the "Metrics" sheet's B2:B13 values (source data refresh).
# All downstream formulas (today!B11:B22, E11:E22, F11:F22, etc. which
# reference Metrics!Bn) recalc automatically.
$wb = $excel.ActiveWorkbook
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value = 92293.209999999992
$wsMetrics.Range("B3").Value = 83491.429999999993
$wsMetrics.Range("B4").Value = 31606.7
$wsMetrics.Range("B5").Value = 3748
$wsMetrics.Range("B6").Value = 672928.98
$wsMetrics.Range("B7").Value = 536285.88
$wsMetrics.Range("B8").Value = 196041.26
$wsMetrics.Range("B9").Value = 27224
$wsMetrics.Range("B10").Value = 34774180.700000003
$wsMetrics.Range("B11").Value = 32582278.669999998
$wsMetrics.Range("B12").Value = 12141855.119999999
$wsMetrics.Range("B13").Value = 1345131

# Update the saved cell-selection cursor on the "Metrics" sheet. Selecting
# a range briefly activates its worksheet (matches Excel COM semantics),
# so the "today" sheet is re-activated/re-selected afterwards to restore
# it as the workbook's active tab.
$wsMetrics.Range("C20").Select()

$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("D9").Select()
